$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(49).Insert()

$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 44607
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = 100112017
$ws.Range("G49").Value = "Apio"
$ws.Range("H49").Value = "Americana (o)"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 60
$ws.Range("K49").Value = 8000
$ws.Range("L49").Value = 8500
$ws.Range("M49").Value = 8250
$ws.Range("N49").Value = "$/docena de matas"
$ws.Range("O49").Value = "Provincia del Elquí"
$ws.Range("P49").Value = 1375
$ws.Range("Q49").Value = 6
$ws.Range("R49").Value = "Hortaliza"
